# Insert a new daily price record for Cebollín (Terminal La Palmera de La
# Serena) before the existing row 232, shifting all following rows down by
# one. This brings the weekly-sampled series to a denser (weekly -> more
# frequent) cadence per the commit message "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 232 (and everything after it) down by one row.
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new record.
$ws.Range("A232").Value = 8
$ws.Range("B232").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C232").Value = 'Coquimbo'
$ws.Range("D232").Value = 44988
$ws.Range("E232").Value = 4
$ws.Range("F232").Value = 100112037
$ws.Range("G232").Value = 'Cebollín'
$ws.Range("H232").Value = 'Sin especificar'
$ws.Range("I232").Value = 'Primera'
$ws.Range("J232").Value = 1200
$ws.Range("K232").Value = 1200
$ws.Range("L232").Value = 1400
$ws.Range("M232").Value = 1300
$ws.Range("N232").Value = '$/paquete 6 unidades'
$ws.Range("O232").Value = 'Provincia del Elquí'
$ws.Range("P232").Value = 217
$ws.Range("Q232").Value = 6
$ws.Range("R232").Value = 'Hortaliza'
